$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(
    124.6713943481445,
    122.3664627075195,
    119.7920227050781,
    121.3308792114258,
    127.6945419311523,
    117.9959106445312,
    125.3270645141602,
    122.6628952026367,
    126.8737945556641,
    144.7876281738281,
    150.7129821777344,
    189.806640625,
    197.7531890869141,
    198.4947662353516,
    175.6856689453125,
    158.7167358398438,
    140.3737335205078,
    125.7741546630859,
    115.9532928466797,
    94.09967041015625,
    96.89784240722656,
    85.75377655029297,
    80.98622131347656,
    68.49890899658203
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
